$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

$ws.Hyperlinks.Add($ws.Range("B2"), "http://example.com")
$ws.Hyperlinks.Item(1).Delete()
$ws.Range("B2").Value = ""

$ws.Range("A2").Value = "'"
$ws.Range("A2").Value = ""

$ws.Range("B2").Select()
